$d = $word.ActiveDocument

# Op 1: replace "We've also added the ability to draw freeform " paragraph (10)
# with the Synthesiability / Familiarity / Susan sections.
$d.Paragraphs(10).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Synthesiability</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>The goal of the interface is to create polygons around objects in the image and the interface provides immediate response to indicate whether or not the user is reaching that goal (it draws the next segment and point whenever the user clicks)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>3. Familiarity:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>There are two types of personas who will generally be using the labeller. One type of person could be Susan, a middle aged woman who has been asked to participate in the research project by her son who works at the University. Susan has had very little experience working with computers and approaches the interface with almost entirely fresh perspective.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">For Susan, concept of drawing will most clearly relate to drawing with a pen and pencil in the real world. The interface should thus make </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>affordances</w:t></w:r><w:r><w:t xml:space="preserve"> for her past knowledge. If given a piece of paper with an image and told to draw polygons around the objects on the page, Susan would </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

# Op 2: merge/re-split the "In our version of the labeller..." paragraph (8)
$d.Paragraphs(8).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">In our version of the labeller, we remove the new object button, instead having the user complete a polygon by clicking the starting point. In doing so we reduce the clickable width to approximately ¼ of the original, but make up for it given the dramatic reduction in distance (remember that the user must travel to the finish object button, and then back again to the image). Therefore, assuming that we start in the center, we’ve reduced traveling distance by about ¼ as well. </w:t></w:r><w:r><w:t>The real gain, however, is</w:t></w:r><w:r><w:t xml:space="preserve"> in the form of task continuity and familiarity (which is discussed later)</w:t></w:r></w:p>')

# Op 3: rewrite the "Since the user will be..." paragraph (4)
$d.Paragraphs(4).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Although this behavior is predictable, s</w:t></w:r><w:r><w:t xml:space="preserve">ince the user will be </w:t></w:r><w:r><w:t>circling the object with a polygon</w:t></w:r><w:r><w:t>, his or her mouse will be pointed towards the image. The mouse will therefore be the entire distance of the image on average away from the button. Fitt’s law demonstrates a negative correlation between the time it takes the user to move the mouse to a required location and the distance the object is from the starting point.</w:t></w:r></w:p>')

# Op 4: strip the trailing "In the original interface, ..." runs from the
# drawing paragraph (3), keeping only the picture run, and reposition the image.
$d.Paragraphs(3).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251658240" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="0426893C" wp14:editId="5698F51C"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>0</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>200025</wp:posOffset></wp:positionV><wp:extent cx="3459480" cy="2352040"/><wp:effectExtent l="0" t="0" r="0" b="10160"/><wp:wrapSquare wrapText="bothSides"/><wp:docPr id="1" name="Picture 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 1"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId5"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="3459480" cy="2352040"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="page"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="page"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></w:r></w:p>')
$d.Shapes(1).Top = 11.2

# Op 5: insert the new "Predictability" section between the intro paragraph (2)
# and the drawing paragraph (paragraph 2 itself is kept empty, unchanged).
$d.Paragraphs(2).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Applying Alan Dix’s principles of learnability to the labeller we can determine how best to implement the interface. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>1. Predictability:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>The provided interface is very predictable. Clicking on a point on the image, always creates a point and draws a line to the previous point. Clicking the “finish polygon” button consistently finishes the polygon.</w:t></w:r></w:p>')

# Op 6: split "steps" -> "tasks" in the opening paragraph (1) into 3 runs.
$d.Paragraphs(1).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The image labeller asks the user to spend time completing a repetitive task for no compensation. Therefore the key focus throughout the development of this user interface is the reduction of </w:t></w:r><w:r><w:t xml:space="preserve">tasks </w:t></w:r><w:r><w:t>needed to identify and label and image.</w:t></w:r></w:p>')

